# Final approval of act recognition review items:
# change Status column (H) for rows 8-11 from "Open" to "Closed".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H8").Value = "Closed"
$ws.Range("H9").Value = "Closed"
$ws.Range("H10").Value = "Closed"
$ws.Range("H11").Value = "Closed"
